# "Correcciones en la presentacion"
#
# 1) Slide 17 ("LINQ") is repurposed in-place to become the new
#    "Metodos de extension" slide.
# 2) A new slide 18 ("Funciones anonimas (Lambda)") is inserted.
# 3) A new slide 19 ("LINQ", expanded with a "Formato" section) is
#    inserted, carrying forward the bullet content that used to live
#    on the old slide 17.

$p = $ppt.ActivePresentation

function Set-WingdingsBullets($textRange) {
    $pf = $textRange.ParagraphFormat
    $pf.Bullet.Visible = $true
    $pf.Bullet.Character = [char]252
    $pf.Bullet.Font.Name = "Wingdings"
}

# ---------------------------------------------------------------
# 1) Slide 17: LINQ -> Metodos de extension
# ---------------------------------------------------------------
$s17 = $p.Slides.Item(17)

$title17 = $s17.Shapes.Item(1).TextFrame.TextRange
$title17.Text = "Metodos"
$title17.InsertAfter(" de ")
$title17.InsertAfter("extension ")

$body17 = $s17.Shapes.Item(2).TextFrame.TextRange
$body17.Text = "Se "
$body17.InsertAfter("incorporan")
$body17.InsertAfter(" ")
$body17.InsertAfter("en")
$body17.InsertAfter(" C# 3.0")

$body17.InsertAfter("`rSe ")
$body17.InsertAfter("utiliza")
$body17.InsertAfter(" la palabra this para ")
$body17.InsertAfter("identificar")
$body17.InsertAfter(" el ")
$body17.InsertAfter("tipo")
$body17.InsertAfter(" que se ")
$body17.InsertAfter("extiende")

$body17.InsertAfter("`rSolo ")
$body17.InsertAfter("pueden")
$body17.InsertAfter(" accede las ")
$body17.InsertAfter("propiedades")
$body17.InsertAfter(" ")
$body17.InsertAfter("publicas")
$body17.InsertAfter(" de ")
$body17.InsertAfter("los")
$body17.InsertAfter(" ")
$body17.InsertAfter("tipos")
$body17.InsertAfter(" que ")
$body17.InsertAfter("extienden")

$body17.InsertAfter("`rLa firma no ")
$body17.InsertAfter("puede")
$body17.InsertAfter(" ")
$body17.InsertAfter("coincidir")
$body17.InsertAfter(" con ")
$body17.InsertAfter("ningun")
$body17.InsertAfter(" ")
$body17.InsertAfter("metodo")
$body17.InsertAfter(" del ")
$body17.InsertAfter("tipo")
$body17.InsertAfter(" que ")
$body17.InsertAfter("extiende")

$body17.InsertAfter("`rTienen")
$body17.InsertAfter(" que ")
$body17.InsertAfter("estar")
$body17.InsertAfter(" ")
$body17.InsertAfter("en")
$body17.InsertAfter(" ")
$body17.InsertAfter("una")
$body17.InsertAfter(" ")
$body17.InsertAfter("clase")
$body17.InsertAfter(" ")
$body17.InsertAfter("estatica")

Set-WingdingsBullets $body17

# ---------------------------------------------------------------
# 2) New slide 18: Funciones anonimas (Lambda)
# ---------------------------------------------------------------
$s18 = $p.Slides.Add(18, 2)

$title18 = $s18.Shapes.Item(1).TextFrame.TextRange
$title18.Text = "Funciones"
$title18.InsertAfter(" ")
$title18.InsertAfter("anonimas")
$title18.InsertAfter(" (Lambda)")

$body18 = $s18.Shapes.Item(2).TextFrame.TextRange
$body18.Text = "Utilizan"
$body18.InsertAfter(" el ")
$body18.InsertAfter("operador")
$body18.InsertAfter(" =>")
$body18.InsertAfter("`rTipos")
$body18.InsertAfter("`rFunc")
$body18.InsertAfter("`rAction")
$body18.Paragraphs(3, 1).IndentLevel = 2
$body18.Paragraphs(4, 1).IndentLevel = 2
Set-WingdingsBullets $body18

# ---------------------------------------------------------------
# 3) New slide 19: LINQ (Formato + Funciones)
# ---------------------------------------------------------------
$s19 = $p.Slides.Add(19, 2)

$title19 = $s19.Shapes.Item(1).TextFrame.TextRange
$title19.Text = "LINQ"

$body19 = $s19.Shapes.Item(2).TextFrame.TextRange
$body19.Text = "Formato" + "`r" + "Lamda" + "`r" + "Query" + "`r" + `
    "Funciones" + "`r" + "Where" + "`r" + "Order" + "`r" + "Join" + "`r" + "Let"
$body19.Paragraphs(2, 1).IndentLevel = 2
$body19.Paragraphs(3, 1).IndentLevel = 2
$body19.Paragraphs(5, 1).IndentLevel = 2
$body19.Paragraphs(6, 1).IndentLevel = 2
$body19.Paragraphs(7, 1).IndentLevel = 2
$body19.Paragraphs(8, 1).IndentLevel = 2
Set-WingdingsBullets $body19
